$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")
$ws2 = $wb.Worksheets.Item("Jason Schema")

# --- Update client names on Weekly Timesheet (column B, rows 2-6) ---
$ws1.Range("B2").Value = "Corr"
$ws1.Range("B3").Value = "Moulton"
$ws1.Range("B4").Value = "Jackson / Ho"
$ws1.Range("B5").Value = "Smith"
$ws1.Range("B6").Value = "McGill"

# --- Zero out Rate (col E) and Total (col F) for rows 2-6 on Weekly Timesheet ---
for ($r = 2; $r -le 6; $r++) {
    $ws1.Cells.Item($r, 5).Value = 0
    $ws1.Cells.Item($r, 6).Value = 0
}

# --- Zero out subtotal/grand-total Total cells on Weekly Timesheet ---
$ws1.Range("F8").Value = 0
$ws1.Range("F11").Value = 0
$ws1.Range("F13").Value = 0

# --- Mirror the same client names / rate / total changes on Jason Schema sheet ---
$ws2.Range("D2").Value = "Corr"
$ws2.Range("D3").Value = "Moulton"
$ws2.Range("D4").Value = "Jackson / Ho"
$ws2.Range("D5").Value = "Smith"
$ws2.Range("D6").Value = "McGill"

for ($r = 2; $r -le 6; $r++) {
    $ws2.Cells.Item($r, 6).Value = 0
    $ws2.Cells.Item($r, 7).Value = 0
    # Clear the "Notes" cell but keep it as a (blank) text cell rather than
    # deleting it outright: write a bare quote-prefix then restore the
    # default "Normal" cell style so no stray formatting is left behind.
    $noteCell = $ws2.Cells.Item($r, 9)
    $noteCell.Value = "'"
    $noteCell.Style = "Normal"
}

# --- Employee ID update on Jason Schema sheet (column B, rows 2-6) ---
for ($r = 2; $r -le 6; $r++) {
    $ws2.Cells.Item($r, 2).Value = "emp_pu67gtu5"
}
